$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 80; this shifts existing rows 80:221 down to 81:222
# and Excel automatically extends the sheet dimension to A1:R222.
$ws.Rows("80:80").Insert()

# Populate the newly inserted row 80 with the new record.
$ws.Cells.Item(80, 1).Value = 5
$ws.Cells.Item(80, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(80, 3).Value = "Maule"
$ws.Cells.Item(80, 4).Value = 44540
$ws.Cells.Item(80, 5).Value = 7
$ws.Cells.Item(80, 6).Value = 100114014
$ws.Cells.Item(80, 7).Value = "Betarraga"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 5000
$ws.Cells.Item(80, 11).Value = 500
$ws.Cells.Item(80, 12).Value = 500
$ws.Cells.Item(80, 13).Value = 500
$ws.Cells.Item(80, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(80, 15).Value = "Región del Maule"
$ws.Cells.Item(80, 16).Value = 100
$ws.Cells.Item(80, 17).Value = 5
$ws.Cells.Item(80, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date-number-format style used by the rest
# of column D (same style index as the cells it pushed down).
$ws.Cells.Item(80, 4).NumberFormat = $ws.Cells.Item(81, 4).NumberFormat
